$wb = $excel.ActiveWorkbook

# --- 1. Rename "TimeSlot" sheet to "EDTSlot" (used for ED time slots) ---
$edt = $wb.Worksheets.Item("TimeSlot")
$edt.Name = "EDTSlot"

# relabel the slot names from generic SLOT# to EDT#
$edt.Cells.Item(2, 2).Value = "EDT1"
$edt.Cells.Item(3, 2).Value = "EDT2"
$edt.Cells.Item(4, 2).Value = "EDT3"
$edt.Cells.Item(5, 2).Value = "EDT4"
$edt.Cells.Item(6, 2).Value = "EDT5"
$edt.Cells.Item(7, 2).Value = "EDT6"

# leave the cursor where the author last left it on this sheet
$edt.Range("H23").Select()

# --- 2. Duplicate the EDTSlot sheet to create the UC time-slot sheet ---
$edt.Copy([System.Reflection.Missing]::Value, $edt)
$uct = $wb.Worksheets.Item($edt.Index + 1)
$uct.Name = "UCTSlot"

# relabel the slot names for UC
$uct.Cells.Item(2, 2).Value = "UCT1"
$uct.Cells.Item(3, 2).Value = "UCT2"
$uct.Cells.Item(4, 2).Value = "UCT3"
$uct.Cells.Item(5, 2).Value = "UCT4"
$uct.Cells.Item(6, 2).Value = "UCT5"
$uct.Cells.Item(7, 2).Value = "UCT6"

# leave the cursor where the author last left it on this sheet, and make
# this the active sheet/tab, matching the final saved workbook state
$uct.Range("O17").Select()
$uct.Activate()

Write-Output "Added EDTSlot/UCTSlot timeslot sheets for UC, ED"
